# Adds speaker notes to each of the 4 slides in the presentation.
# (Matches the commit that introduced ppt/notesSlides/notesSlide1-4.xml,
#  describing the Java Swing component hierarchy used for each mock-up.)

$p = $ppt.ActivePresentation

$notesText1 = "JFrame`n`tJTabbedPane`n`t`tJScrollPane`n`t`t`tJTable`n`tJButton (“Create New Chat”)`n`tJButton (“Log Out”)`n"
$notesText2 = "JFrame`n`tJTabbedPane`n`t`tJScrollPane`n`t`t`tJTree`n`tJButton (“Create New Chat”)`n`tJButton (“Log Out”)`n"
$notesText3 = "JFrame`n`tJTabbedPane`n`t`tJScrollPane`n`t`t`tJTable`n`t`t`t`tJButton`n`t`t`t`tJButton`n`tJButton (“Create New Chat”)`n`tJButton (“Log Out”)`n"
$notesText4 = "JFrame`n`tJTabbedPane`n`t`tJScrollPane`n`t`t`tJTable`n`tJButton (“Create New Chat”)`n`tJButton (“Log Out”)"

$notesBySlide = @{ 1 = $notesText1; 2 = $notesText2; 3 = $notesText3; 4 = $notesText4 }

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $notesPage = $slide.NotesPage
    # Only the notes-body placeholder is writable on a notes page; this both
    # materializes the notes slide part and sets its text in one call.
    $notesShape = $notesPage.Shapes.AddPlaceholder(2)
    $notesShape.TextFrame.TextRange.Text = $notesBySlide[$i]
}
